$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at the top, shifting existing data (rows 1-30) down to rows 2-31.
$ws.Rows.Item(1).Insert()

# Populate the new header row with "Tasks" / "Weights".
$ws.Range("A1").Value = "Tasks"
$ws.Range("B1").Value = "Weights"

# Clear the style that Insert() may have copied into B1 from the old row 1 (B2),
# so the new header cell uses the default (unstyled) formatting.
$ws.Range("B1").Style = "Normal"

# Reset the view: scroll back to the top and select B1.
$ws.Application.ActiveWindow.ScrollRow = 1
$ws.Range("B1").Select()
